# Feature/22 improve dealloc policy (#23)
# - Merge previously-split text runs back into single runs on a few slides
#   (no visible text change, just run consolidation as produced by PowerPoint
#   when text is retyped/edited in place).
# - Add a new slide (hint on addr2line and nm usage).

$p = $ppt.ActivePresentation

function Merge-Runs($Shape, $OldText) {
    $tr = $Shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($OldText)
    if ($idx -lt 0) {
        throw "Merge-Runs: text not found: $OldText"
    }
    $sub = $tr.Characters($idx + 1, $OldText.Length)
    $sub.Text = $OldText
}

# --- slide 2 ("Overview"): merge "Memory " + "handling related " + "crashes" ---
$s2 = $p.Slides.Item(2)
Merge-Runs $s2.Shapes.Item(2) "Memory handling related crashes"

# --- slide 4 ("How Valgrind works..."): merge "of " + "concurrency." ---
$s4 = $p.Slides.Item(4)
Merge-Runs $s4.Shapes.Item(2) "of concurrency."

# --- slide 5 ("Crash investigator - idea behind") ---
$s5 = $p.Slides.Item(5)
Merge-Runs $s5.Shapes.Item(1) "Crash investigator - idea behind"
Merge-Runs $s5.Shapes.Item(2) "In order to trap memory allocations functions and make some analyze the following can be done"

# --- slide 6 ("Some cases of indirect double/free") ---
$s6 = $p.Slides.Item(6)
Merge-Runs $s6.Shapes.Item(2) "Global buffers allocation/deallocation without proper synchronization."

# --- add new slide 8: addr2line / nm usage hint ---
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)
$newSlide.Shapes.Item(1).Name = "Title 1"
$newSlide.Shapes.Item(2).Name = "Content Placeholder 2"

$titleTr = $newSlide.Shapes.Item(1).TextFrame.TextRange
$titleTr.Text = ""

$bodyTr = $newSlide.Shapes.Item(2).TextFrame.TextRange
$bodyTr.Text = "addr2line -e /home/kalantar/dev/crash_investigator/sys/focal/Debug/lib/libcrash_investigator_new_malloc_0020.so -f -C "
$bodyTr.LanguageID = "de-DE"
$r2 = $bodyTr.InsertAfter("0x7d75")
$r2.LanguageID = "de-DE"
$r3 = $r2.InsertAfter("`n")
$r3.LanguageID = "de-DE"
$r4 = $r3.InsertAfter("nm")
$r4.LanguageID = "de-DE"
$r5 = $r4.InsertAfter(" /home/kalantar/dev/crash_investigator/sys/focal/Debug/lib/libcrash_investigator_new_malloc_0020.so | ")
$r5.LanguageID = "de-DE"
$r6 = $r5.InsertAfter("grep")
$r6.LanguageID = "de-DE"
$r7 = $r6.InsertAfter(" _ZN18crash_investigator11CMemoryItem4InitEmNS_11FailureTypeEPvPNS_9BacktraceE")
$r7.LanguageID = "de-DE"
